$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Read the template values from row 1 (column A-G) using .Text, since
# .Value getters return an unresolved placeholder in this environment.
$colValues = @{}
for ($col = 1; $col -le 7; $col++) {
    $colValues[$col] = $ws.Cells.Item(1, $col).Text
}

# Duplicate row 1's A:G content into new rows 3, 4 and 5, and add a
# newline-only value in column H for each of those rows.
for ($r = 3; $r -le 5; $r++) {
    for ($col = 1; $col -le 7; $col++) {
        $ws.Cells.Item($r, $col).Value = $colValues[$col]
    }
    $ws.Cells.Item($r, 8).Value = "`n"
}
